# Add: knn using sklearn
#
# Turn the raw Hue/Saturation/Value table (A1:C17, no header) into a
# labelled dataset for a KNN classifier:
#   - add a header row: Hue, Saturation, Value, Class
#   - drop 3 stray rows (old rows 8, 11, 17)
#   - add a Class column (D) labelling the first 7 remaining rows "apples"
#     and the rest "orange"
# Final extent: A1:D15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the existing raw values (17 rows x 3 cols, no header).
# Range.Value() returns a 1-based 2D array: old[r,1..3].
$old = $ws.Range("A1:C17").Value()

# Rows (1-based, in the ORIGINAL sheet) to drop.
$dropRows = @(8, 11, 17)

$keptRows = @()
for ($r = 1; $r -le 17; $r++) {
    if ($dropRows -notcontains $r) {
        $keptRows += $r
    }
}

# Clear the old footprint (A1:C17) since the new one (A1:D15) is smaller.
$ws.Range("A1:C17").Clear()

$ws.Cells.Item(1, 1).Value = "Hue"
$ws.Cells.Item(1, 2).Value = "Saturation"
$ws.Cells.Item(1, 3).Value = "Value"
$ws.Cells.Item(1, 4).Value = "Class"

$destRow = 2
foreach ($srcRow in $keptRows) {
    $ws.Cells.Item($destRow, 1).Value = $old[$srcRow, 1]
    $ws.Cells.Item($destRow, 2).Value = $old[$srcRow, 2]
    $ws.Cells.Item($destRow, 3).Value = $old[$srcRow, 3]
    if ($destRow -le 8) {
        $ws.Cells.Item($destRow, 4).Value = "apples"
    } else {
        $ws.Cells.Item($destRow, 4).Value = "orange"
    }
    $destRow++
}
